$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.552.17"
$ws.Range("E2").Value = "  -0.78%  "
$ws.Range("D3").Value = "1.859.34"
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("E4").Value = "  +0.79%  "
$ws.Range("D5").Value = "333.14"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "1.010"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").Value = "0.4660"
$ws.Range("E7").Value = "  -1.20%  "
$ws.Range("D8").Value = "0.3883"
$ws.Range("E8").Value = "  -1.88%  "
$ws.Range("D9").Value = "45.74"
$ws.Range("E9").Value = "  -4.45%  "
$ws.Range("D10").Value = "0.07951"
$ws.Range("E10").Value = "  -1.28%  "
$ws.Range("E11").Value = "  -3.32%  "
$ws.Range("D12").Value = "21.57"
$ws.Range("E12").Value = "  -2.90%  "
$ws.Range("D13").Value = "1.863.14"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "5.973"
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").Value = "7.197"
$ws.Range("E15").Value = "  +0.96%  "
$ws.Range("D16").Value = "1.014"
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("D17").Value = "87.68"
$ws.Range("E17").Value = "  +0.48%  "
$ws.Range("D18").Value = "0.06707"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("D19").Value = "0.00001040"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").Value = "16.88"
$ws.Range("E20").Value = "  -2.52%  "
$ws.Range("D21").Value = "1.010"
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("D22").Value = "27.538.90"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").Value = "5.428"
$ws.Range("E23").Value = "  -1.95%  "
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("D25").Value = "2.302"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").Value = "2.083.87"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").Value = "158.76"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("D28").Value = "19.68"
$ws.Range("E28").Value = "  -2.51%  "
$ws.Range("D29").Value = "2.118"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("D30").Value = "5.371"
$ws.Range("E30").Value = "  -3.86%  "
$ws.Range("D31").Value = "120.99"
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("D32").Value = "0.9697"
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("D33").Value = "0.09442"
$ws.Range("E33").Value = "  -0.96%  "
$ws.Range("D34").Value = "3.641"
$ws.Range("E34").Value = "  +1.20%  "
$ws.Range("D35").Value = "5.282"
$ws.Range("E35").Value = "  -1.50%  "
$ws.Range("D36").Value = "1.323"
$ws.Range("E36").Value = "  -8.76%  "
$ws.Range("D37").Value = "0.06018"
$ws.Range("E37").Value = "  -1.90%  "
$ws.Range("D38").Value = "0.02209"
$ws.Range("E38").Value = "  -2.21%  "
$ws.Range("D39").Value = "1.189"
$ws.Range("E39").Value = "  -3.81%  "
$ws.Range("D40").Value = "8.179"
$ws.Range("E40").Value = "  +0.37%  "
$ws.Range("D41").Value = "1.010"
$ws.Range("E41").Value = "  +0.76%  "
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("D43").Value = "0.1873"
$ws.Range("E43").Value = "  -1.55%  "
$ws.Range("D44").Value = "10.17"
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("D45").Value = "1.240"
$ws.Range("E45").Value = "  -2.60%  "
$ws.Range("D46").Value = "0.5602"
$ws.Range("E46").Value = "  -2.00%  "
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("D48").Value = "1.907"
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("D49").Value = "3.263"
$ws.Range("E49").Value = "  -3.48%  "
$ws.Range("D50").Value = "0.06753"
$ws.Range("E50").Value = "  -2.43%  "
$ws.Range("D51").Value = "112.16"
$ws.Range("E51").Value = "  -1.94%  "
